$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(29, 58, 99, 127, 163, 213, 251, 298, 341)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$wb.Save()
